$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "45.426.14"
Set-TextValue "E2" "  +7.26%  "
Set-TextValue "D3" "2.391.04"
Set-TextValue "E3" "  +4.84%  "
Set-TextValue "E4" "  +0.99%  "
Set-TextValue "D5" "113.25"
Set-TextValue "E5" "  +9.63%  "
Set-TextValue "D6" "318.47"
Set-TextValue "E6" "  +2.63%  "
Set-TextValue "D7" "0.638"
Set-TextValue "E7" "  +2.85%  "
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D9" "0.629"
Set-TextValue "E9" "  +4.98%  "
Set-TextValue "D10" "42.51"
Set-TextValue "E10" "  +10.35%  "
Set-TextValue "D11" "0.0934"
Set-TextValue "E11" "  +4.03%  "
Set-TextValue "E12" "  +6.15%  "
Set-TextValue "E13" "  +5.59%  "
Set-TextValue "E14" "  +1.82%  "
Set-TextValue "D15" "15.90"
Set-TextValue "E15" "  +5.71%  "
Set-TextValue "D16" "2.753.24"
Set-TextValue "E16" "  +4.79%  "
Set-TextValue "D17" "2.390.13"
Set-TextValue "E17" "  +4.72%  "
Set-TextValue "D18" "45.415.68"
Set-TextValue "E18" "  +7.30%  "
Set-TextValue "E19" "  +6.01%  "
Set-TextValue "E20" "  +4.15%  "
Set-TextValue "D21" "13.44"
Set-TextValue "E21" "  +3.98%  "
Set-TextValue "D22" "75.09"
Set-TextValue "E22" "  +3.33%  "
Set-TextValue "E23" "  +5.12%  "
Set-TextValue "D24" "269.98"
Set-TextValue "E24" "  +2.97%  "
Set-TextValue "D25" "2.38"
Set-TextValue "E25" "  +9.50%  "
Set-TextValue "E26" "  -0.72%  "
Set-TextValue "D27" "11.29"
Set-TextValue "E27" "  +6.24%  "
Set-TextValue "E28" "  +9.24%  "
Set-TextValue "E29" "  +2.30%  "
Set-TextValue "D30" "39.07"
Set-TextValue "E30" "  +9.68%  "
Set-TextValue "D31" "22.94"
Set-TextValue "E31" "  +3.92%  "
Set-TextValue "D32" "0.0943"
Set-TextValue "E32" "  +10.88%  "
Set-TextValue "D33" "170.54"
Set-TextValue "E33" "  +3.67%  "
Set-TextValue "D34" "2.97"
Set-TextValue "E34" "  +16.88%  "
Set-TextValue "E35" "  +3.80%  "
Set-TextValue "D36" "4.92"
Set-TextValue "E36" "  +9.98%  "
Set-TextValue "D37" "0.118"
Set-TextValue "E37" "  +7.54%  "
Set-TextValue "D38" "3.15"
Set-TextValue "E38" "  +17.28%  "
Set-TextValue "D39" "0.0365"
Set-TextValue "E39" "  +5.44%  "
Set-TextValue "E40" "  +8.16%  "
Set-TextValue "E41" "  +13.48%  "
Set-TextValue "D42" "104.92"
Set-TextValue "E42" "  +6.78%  "
Set-TextValue "D43" "0.242"
Set-TextValue "E43" "  +7.79%  "
Set-TextValue "D44" "13.58"
Set-TextValue "E44" "  +14.51%  "
Set-TextValue "D45" "71.91"
Set-TextValue "E45" "  +5.11%  "
Set-TextValue "E46" "  +0.12%  "
Set-TextValue "D47" "117.86"
Set-TextValue "E47" "  +7.41%  "
Set-TextValue "E48" "  +14.03%  "
Set-TextValue "E49" "  +21.24%  "
Set-TextValue "D50" "9.39"
Set-TextValue "E50" "  +9.30%  "
Set-TextValue "D51" "0.226"
Set-TextValue "E51" "  +19.73%  "
